# Update column G ("K") values for rows 2-16 on the active sheet.
# These are literal data values (strikeouts recalculated), not formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 2
    3  = 1
    4  = 0
    5  = 4
    6  = 3
    7  = 1
    8  = 1
    9  = 2
    10 = 2
    12 = 0
    13 = 1
    14 = 2
    15 = 2
    16 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
